$wb = $excel.ActiveWorkbook

$wsFicha = $wb.Worksheets.Item("Ficha técnica")

$wsFicha.Range("B5").Value = "El indicador mide el porcentaje de niños, niñas y adolescentes (hasta 17 años de edad) que residen en hogares en los que ningún integrante percibe Asignaciones Familaires - Plan de Equidad ni Tarjeta Uruguay Social en el total de hogares pobres"

$wsFicha.Range("B6").Value = "Para cada año calcular: (Cantidad de niños, niñas y adolescentes (hasta 17 años de edad) que residen en hogares en los que ninguno de sus integrantes es beneficiario de Asignaciones Familiares - Plan de Equidad ni Tarjeta Uruguay Social/Cantidad de niños, niñas y adolescentes que residen en hogares pobres)*100"
